# Weekly update: add this week's Primera/Segunda price records for
# Vega Monumental Concepción - Arándano (blue).
#
# Two new rows are inserted right before the current row 49, pushing all
# existing data rows (and the sheet dimension) down by two. The new rows
# are then populated with the new weekly observations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 49:50 (existing rows 49.. shift down to 51..).
$ws.Range("A49:A50").EntireRow.Insert()

# --- New row 49: "Primera" quality ---
$ws.Cells.Item(49, 1).Value = 11
$ws.Cells.Item(49, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(49, 3).Value = "Bíobío"
$ws.Cells.Item(49, 4).Value = 44574
$ws.Cells.Item(49, 5).Value = 8
$ws.Cells.Item(49, 6).Value = "Fruta"
$ws.Cells.Item(49, 7).Value = 100101
$ws.Cells.Item(49, 8).Value = "Berries"
$ws.Cells.Item(49, 9).Value = 100101001
$ws.Cells.Item(49, 10).Value = "Arándano (blue)"
$ws.Cells.Item(49, 11).Value = "Sin especificar"
$ws.Cells.Item(49, 12).Value = "Primera"
$ws.Cells.Item(49, 13).Value = 200
$ws.Cells.Item(49, 14).Value = 3500
$ws.Cells.Item(49, 15).Value = 4000
$ws.Cells.Item(49, 16).Value = 3750
$ws.Cells.Item(49, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(49, 18).Value = "Región de Ñuble"
$ws.Cells.Item(49, 19).Value = 1875
$ws.Cells.Item(49, 20).Value = 2

# --- New row 50: "Segunda" quality ---
$ws.Cells.Item(50, 1).Value = 11
$ws.Cells.Item(50, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(50, 3).Value = "Bíobío"
$ws.Cells.Item(50, 4).Value = 44574
$ws.Cells.Item(50, 5).Value = 8
$ws.Cells.Item(50, 6).Value = "Fruta"
$ws.Cells.Item(50, 7).Value = 100101
$ws.Cells.Item(50, 8).Value = "Berries"
$ws.Cells.Item(50, 9).Value = 100101001
$ws.Cells.Item(50, 10).Value = "Arándano (blue)"
$ws.Cells.Item(50, 11).Value = "Sin especificar"
$ws.Cells.Item(50, 12).Value = "Segunda"
$ws.Cells.Item(50, 13).Value = 100
$ws.Cells.Item(50, 14).Value = 3000
$ws.Cells.Item(50, 15).Value = 3000
$ws.Cells.Item(50, 16).Value = 3000
$ws.Cells.Item(50, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(50, 18).Value = "Región de Ñuble"
$ws.Cells.Item(50, 19).Value = 1500
$ws.Cells.Item(50, 20).Value = 2

# Make sure the date cells keep the workbook's date number format (style
# index 2 in styles.xml), matching the rest of column D.
$ws.Cells.Item(49, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(50, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
